# feat: add 2022-Q1 data
#
# Before: Worksheets = [ "2021-Q4"(id1), "总计"(id2) ]
# After : Worksheets = [ "2021-Q4"(id1), "2022-Q1"(id2), "总计"(id3) ]
#
# The existing "总计" worksheet is duplicated in place (Worksheets.Copy puts
# the clone immediately after it and gives it a fresh sheetId) so both the
# new "2022-Q1" sheet and the refreshed "总计" sheet start from an identical,
# fully-formatted template (same sheetPr/outlinePr/pageMargins/etc. as the
# original) instead of a generic blank Worksheets.Add() sheet. The original
# keeps sheetId 2 and is renamed to "2022-Q1"; the clone gets the next
# sheetId (3) and is renamed to "总计" - matching the target numbering.
# Both are then cleared and repopulated with their respective tables.
#
# Helper notes:
#  - A plain `Range.Value = "006235"` / `"35.36"` gets auto-coerced to a
#    Number by the COM layer (same as real Excel typing), which would lose
#    leading zeros / change the stored type away from Text. Prefixing the
#    literal with a single quote forces Text storage.
#  - That quote-prefix leaves a `quotePrefix` flag baked into the cell's
#    style. `PasteSpecial` with formats-only (xlPasteFormats = -4122),
#    pulling from an already-correctly-styled donor cell, overwrites the
#    cell's style wholesale (including clearing quotePrefix) without
#    touching the value that was just written - so it's used after every
#    text assignment to land on the exact donor style (header/index style,
#    or the plain unstyled look used by the rest of the data cells).

$xlPasteFormats = -4122

function Set-TextCell($range, [string]$text, $styleDonor) {
    $range.Value = "'" + $text
    $styleDonor.Copy() | Out-Null
    $range.PasteSpecial($xlPasteFormats)
}

function Set-IndexCell($range, $value, $styleDonor) {
    $styleDonor.Copy() | Out-Null
    $range.PasteSpecial($xlPasteFormats)
    $range.Value = $value
}

$wb = $excel.ActiveWorkbook

$wsOld2021 = $wb.Worksheets.Item(1)          # "2021-Q4" - untouched, also a format donor
$donorHeaderStyle = $wsOld2021.Range("B1")   # s=2 (bold / centered / bordered)
$donorPlainText   = $wsOld2021.Range("B2")   # default style, text

# ------------------------------------------------------------------------
# 0) Duplicate "总计" in place so both resulting sheets inherit its exact
#    sheetPr / sheetFormatPr / pageMargins layout.
# ------------------------------------------------------------------------
$wsTotalOrig = $wb.Worksheets.Item(2)        # "总计", sheetId 2
$wsTotalOrig.Copy($null, $wsTotalOrig)       # clone lands right after -> sheetId 3

$wsQ1 = $wb.Worksheets.Item(2)               # the original -> becomes "2022-Q1"
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Item(3)            # the clone -> stays "总计"
$wsTotal.Name = "总计"

# ------------------------------------------------------------------------
# 1) Populate "2022-Q1" with the fund-holding table.
# ------------------------------------------------------------------------
$wsQ1.Cells.Clear()

Set-TextCell $wsQ1.Range("B1") "基金代码"       $donorHeaderStyle
Set-TextCell $wsQ1.Range("C1") "基金名称"       $donorHeaderStyle
Set-TextCell $wsQ1.Range("D1") "基金规模"       $donorHeaderStyle
Set-TextCell $wsQ1.Range("E1") "股票总仓位"     $donorHeaderStyle
Set-TextCell $wsQ1.Range("F1") "仓位占比"       $donorHeaderStyle
Set-TextCell $wsQ1.Range("G1") "持有市值(亿元)" $donorHeaderStyle
Set-TextCell $wsQ1.Range("H1") "仓位排名"       $donorHeaderStyle

Set-IndexCell $wsQ1.Range("A2") 0 $donorHeaderStyle
Set-TextCell  $wsQ1.Range("B2") "161040"                     $donorPlainText
Set-TextCell  $wsQ1.Range("C2") "富国创业板两年定期开放混合" $donorPlainText
Set-TextCell  $wsQ1.Range("D2") "35.36"                      $donorPlainText
Set-TextCell  $wsQ1.Range("E2") "83.63"                      $donorPlainText
Set-TextCell  $wsQ1.Range("F2") "3.54"                       $donorPlainText
Set-TextCell  $wsQ1.Range("G2") "1.2517"                     $donorPlainText
$wsQ1.Range("H2").Value = 7

Set-IndexCell $wsQ1.Range("A3") 1 $donorHeaderStyle
Set-TextCell  $wsQ1.Range("B3") "006235"               $donorPlainText
Set-TextCell  $wsQ1.Range("C3") "东方城镇消费主题混合" $donorPlainText
Set-TextCell  $wsQ1.Range("D3") "0.50"                 $donorPlainText
Set-TextCell  $wsQ1.Range("E3") "90.32"                $donorPlainText
Set-TextCell  $wsQ1.Range("F3") "4.58"                 $donorPlainText
Set-TextCell  $wsQ1.Range("G3") "0.0229"               $donorPlainText
$wsQ1.Range("H3").Value = 8

Set-IndexCell $wsQ1.Range("A4") 2 $donorHeaderStyle
Set-TextCell  $wsQ1.Range("B4") "003279"                         $donorPlainText
Set-TextCell  $wsQ1.Range("C4") "融通沪港深智慧生活灵活配置混合" $donorPlainText
Set-TextCell  $wsQ1.Range("D4") "0.13"                           $donorPlainText
Set-TextCell  $wsQ1.Range("E4") "67.65"                          $donorPlainText
Set-TextCell  $wsQ1.Range("F4") "2.81"                           $donorPlainText
Set-TextCell  $wsQ1.Range("G4") "0.0037"                         $donorPlainText
$wsQ1.Range("H4").Value = 6

# ------------------------------------------------------------------------
# 2) Populate "总计" with the refreshed totals table (new 2022-Q1 row on
#    top, old 2021-Q4 row pushed down).
# ------------------------------------------------------------------------
$wsTotal.Cells.Clear()

Set-TextCell $wsTotal.Range("B1") "日期"           $donorHeaderStyle
Set-TextCell $wsTotal.Range("C1") "持有数量(只)"   $donorHeaderStyle
Set-TextCell $wsTotal.Range("D1") "持有市值(亿元)" $donorHeaderStyle

Set-IndexCell $wsTotal.Range("A2") 0 $donorHeaderStyle
Set-TextCell  $wsTotal.Range("B2") "2022-Q1" $donorPlainText
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 1.28

Set-IndexCell $wsTotal.Range("A3") 1 $donorHeaderStyle
Set-TextCell  $wsTotal.Range("B3") "2021-Q4" $donorPlainText
$wsTotal.Range("C3").Value = 5
$wsTotal.Range("D3").Value = 2.88

# ------------------------------------------------------------------------
# 3) Leave the view the way the original workbook had it: first sheet
#    active, selection parked at A1.
# ------------------------------------------------------------------------
$wsOld2021.Activate()
$wsOld2021.Range("A1").Select() | Out-Null
